# Committee workbook update:
#  - Fix typo in affiliation for Agnieszka Kubik-Komar (row 18):
#      "University of Life Sciences in Lubli" -> "University of Life Sciences in Lublin"
#  - Add a new committee member "Lilla Di Scala" (Johnson & Johnson) as a new row,
#    inserted right after Susanne Strohmaier (old row 21), pushing the local-committee
#    members (Laurence Giullier ... Marcel Wolbers) down by one row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the affiliation typo for Agnieszka Kubik-Komar
$ws.Range("C18").Value = "University of Life Sciences in Lublin"

# Insert a new row for the new committee member, shifting rows 22-30 down to 23-31
$ws.Rows.Item(22).Insert()

$ws.Range("A22").Value = "Lilla"
$ws.Range("B22").Value = "Di Scala"
$ws.Range("C22").Value = "Johnson & Johnson"

# Match the author's last cursor position
$ws.Range("I14").Select()
